$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'24.954.50"
$ws.Range("E2").Value = "'  +2.09%  "

# Row 3
$ws.Range("D3").Value = "'1.700.38"
$ws.Range("E3").Value = "'  +0.86%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.21%  "

# Row 5
$ws.Range("D5").Value = "'315.56"
$ws.Range("E5").Value = "'  -0.02%  "

# Row 6
$ws.Range("E6").Value = "'  +0.22%  "

# Row 7
$ws.Range("D7").Value = "'0.3979"
$ws.Range("E7").Value = "'  +1.78%  "

# Row 8
$ws.Range("D8").Value = "'0.4023"
$ws.Range("E8").Value = "'  -0.17%  "

# Row 9
$ws.Range("D9").Value = "'1.465"
$ws.Range("E9").Value = "'  -1.36%  "

# Row 10
$ws.Range("D10").Value = "'52.99"
$ws.Range("E10").Value = "'  +1.21%  "

# Row 11
$ws.Range("D11").Value = "'1.003"
$ws.Range("E11").Value = "'  +0.27%  "

# Row 12
$ws.Range("D12").Value = "'0.08799"
$ws.Range("E12").Value = "'  +0.32%  "

# Row 13
$ws.Range("D13").Value = "'25.99"
$ws.Range("E13").Value = "'  -1.94%  "

# Row 14
$ws.Range("D14").Value = "'7.462"
$ws.Range("E14").Value = "'  +0.17%  "

# Row 15
$ws.Range("D15").Value = "'0.00001352"

# Row 16
$ws.Range("D16").Value = "'7.948"
$ws.Range("E16").Value = "'  -2.48%  "

# Row 17
$ws.Range("D17").Value = "'1.707.83"
$ws.Range("E17").Value = "'  +1.77%  "

# Row 18
$ws.Range("D18").Value = "'95.99"
$ws.Range("E18").Value = "'  -2.16%  "

# Row 19
$ws.Range("D19").Value = "'0.07190"
$ws.Range("E19").Value = "'  -0.56%  "

# Row 20
$ws.Range("D20").Value = "'20.63"
$ws.Range("E20").Value = "'  +2.01%  "

# Row 21
$ws.Range("D21").Value = "'7.328"
$ws.Range("E21").Value = "'  +0.72%  "

# Row 22
$ws.Range("E22").Value = "'  +0.17%  "

# Row 23
$ws.Range("D23").Value = "'14.37"
$ws.Range("E23").Value = "'  +0.69%  "

# Row 24
$ws.Range("D24").Value = "'24.961.20"
$ws.Range("E24").Value = "'  +2.13%  "

# Row 25
$ws.Range("D25").Value = "'2.368"
$ws.Range("E25").Value = "'  +1.35%  "

# Row 26
$ws.Range("D26").Value = "'2.941"
$ws.Range("E26").Value = "'  -3.00%  "

# Row 27
$ws.Range("D27").Value = "'23.76"
$ws.Range("E27").Value = "'  +5.05%  "

# Row 28
$ws.Range("D28").Value = "'6.204"
$ws.Range("E28").Value = "'  +16.12%  "

# Row 29
$ws.Range("D29").Value = "'162.08"
$ws.Range("E29").Value = "'  -3.31%  "

# Row 30
$ws.Range("D30").Value = "'150.67"
$ws.Range("E30").Value = "'  +8.98%  "

# Row 31
$ws.Range("D31").Value = "'8.361"
$ws.Range("E31").Value = "'  -1.00%  "

# Row 32
$ws.Range("D32").Value = "'2.624"
$ws.Range("E32").Value = "'  +25.67%  "

# Row 33
$ws.Range("D33").Value = "'1.896.78"
$ws.Range("E33").Value = "'  +1.72%  "

# Row 34
$ws.Range("D34").Value = "'0.08572"
$ws.Range("E34").Value = "'  -2.06%  "

# Row 35
$ws.Range("D35").Value = "'0.03135"
$ws.Range("E35").Value = "'  +4.00%  "

# Row 36
$ws.Range("B36").Value = "'InternetComputer(DFINITY)"
$ws.Range("C36").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'7.148"
$ws.Range("E36").Value = "'  -1.73%  "

# Row 37
$ws.Range("B37").Value = "'ImmutableX"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.036"
$ws.Range("E37").Value = "'  -1.10%  "

# Row 38
$ws.Range("D38").Value = "'0.2867"
$ws.Range("E38").Value = "'  +2.73%  "

# Row 39
$ws.Range("D39").Value = "'0.09580"
$ws.Range("E39").Value = "'  +4.89%  "

# Row 40
$ws.Range("E40").Value = "'  +0.35%  "

# Row 41
$ws.Range("D41").Value = "'0.8242"

# Row 42
$ws.Range("D42").Value = "'14.01"

# Row 43
$ws.Range("D43").Value = "'1.486"
$ws.Range("E43").Value = "'  +1.10%  "

# Row 44
$ws.Range("D44").Value = "'17.25"
$ws.Range("E44").Value = "'  -1.89%  "

# Row 45
$ws.Range("D45").Value = "'2.683"
$ws.Range("E45").Value = "'  +1.02%  "

# Row 46
$ws.Range("D46").Value = "'0.7382"
$ws.Range("E46").Value = "'  +1.83%  "

# Row 47
$ws.Range("D47").Value = "'4.246"
$ws.Range("E47").Value = "'  -0.38%  "

# Row 48
$ws.Range("D48").Value = "'1.392"
$ws.Range("E48").Value = "'  -1.69%  "

# Row 49
$ws.Range("D49").Value = "'0.08800"
$ws.Range("E49").Value = "'  +8.90%  "

# Row 50
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "'  +0.20%  "

# Row 51
$ws.Range("D51").Value = "'139.24"
$ws.Range("E51").Value = "'  -0.06%  "
